$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.99%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.56%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.15%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07544"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.60%"

# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.604"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.59%"

# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9366"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.33%"

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.384"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.37%"

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1190"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.11%"

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1823"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.47%"

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09075"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.77%"

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04134"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.68%"

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1047"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.77%"

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.64%"

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005839"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.27%"

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.341"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.27%"

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.359"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.65%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3334"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.86%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.329"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "6.56%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1408"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.37%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3097"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "10.60%"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.20%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003900"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "6.08%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001300"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.28%"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02419"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.27%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05218"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.20%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006677"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "11.92%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007710"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.02%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1328"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.82%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007380"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.03%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007104"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.83%"

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.97%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006224"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.10%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03481"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-26.71%"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
